$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two row reorders
# within the Hedera/WEMIXToken/THORChain/LidoDAOToken block and the
# MXToken/RocketPoolETH block) as produced by the latest GitHub Actions run.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '38.254.65'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.068.05'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.04'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.32'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +9.92%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0809'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.74%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.00'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.374.73'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.44'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +7.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.771'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.32'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.74%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '38.215.49'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.33'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.96%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.71'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.19'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.12%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.02'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.33'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.119'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.58'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.69'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.11%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.05'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +9.97%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0611'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'THORChain'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.28'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +15.68%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.29'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.35'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.13%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.526.84'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.33'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.86'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.55%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.45%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0927'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.17'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.44%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.261.43'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.70%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '44.49'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.28%  '
